$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: date (A4) now written as an integer-index date value instead of a
# shared-string timestamp, with the remaining greeks recomputed.
$ws.Range("A4").Value = 43467
$ws.Range("A4").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row4 = @(
    -0.03056031644949887,
    0.009336083822868596,
    0.03989640027236746,
    0.007339653571375165,
    0.003696668535151762,
    -0.003642985036223402,
    -0.02324686710860826,
    -0.0116229998037474,
    0.01162386730486087,
    -0.002351533235689906,
    0.01895979495127779,
    0.0213113281869677
)
$cols = @("B","C","D","E","F","G","H","I","J","K","L","M")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "4").Value = $row4[$i]
}

# Row 5: same treatment.
$ws.Range("A5").Value = 43468
$ws.Range("A5").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$row5 = @(
    -0.01437140681067996,
    0.006029594950748981,
    0.02040100176142894,
    -0.002347481293098961,
    -0.001114083980204146,
    0.001233397312894815,
    -0.004463858392987965,
    -0.002162501257339382,
    0.002301357135648583,
    -0.005054024448394552,
    0.02374142432484177,
    0.02879544877323632
)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "5").Value = $row5[$i]
}
